$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").Value = "Mes (código)"
$ws.Range("B1").Value = "Número de contratos"
$ws.Range("C1").Value = "Provincia código"
$ws.Range("D1").Value = "Provincia nombre"
$ws.Range("E1").Value = "Sexo"
$ws.Range("F1").Value = "Mes y año"

# Row 2
$ws.Range("A2").Value = "null"
$ws.Range("B2").Value = "iaest-measure:numero-de-contratos"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "iaest-measure:sexo"
$ws.Range("F2").Value = "iaest-measure:mes-y-ano"

# Row 3
$ws.Range("A3").Value = "null"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "medida"

# Row 4
$ws.Range("A4").Value = "null"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "URI-Provincia"
$ws.Range("E4").Value = "xsd:string"
$ws.Range("F4").Value = "xsd:string"
